$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.571.75"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "'1.882.95"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'246.39"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.4733"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.2893"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "'0.06540"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'22.38"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'100.05"
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("D12").Value = "'0.7632"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").Value = "'0.07832"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'1.882.65"
$ws.Range("D15").Value = "'5.243"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'284.57"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "'30.543.73"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'0.000007526"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'2.127.39"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "'5.358"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'6.441"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "'9.180"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").Value = "'163.74"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'19.07"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "'1.906"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "'0.09712"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").Value = "'1.329"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'1.500"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "'4.252"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "'4.176"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "'0.04846"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.6979"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "'2.764"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'0.01908"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "'2.872"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").Value = "'6.310"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("D41").Value = "'75.52"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").Value = "'1.980"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "'0.4255"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'0.8402"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "'101.24"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "'9.864"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").Value = "'7.030"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'35.31"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "'0.05780"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'0.3959"
$ws.Range("E51").Value = "  -0.19%  "